$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: handback status text changed from "Ready for handoff"
#    to "Handed back: in sync with en-US" for both zh-cn (E) and de-de (F)
#    columns, on both data rows. Also those two columns get wider.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: the handback finished earlier, but the report generation
#    needs to refresh the handback datetime and clear the stale "not the
#    latest version" error message now that everything is in sync.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("K2").Value = "2016-09-06 15:18:32"
$zhcn.Range("K3").Value = "2016-09-06 15:18:32"

$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------------
# 3. de-de sheet: this is the handback that just completed. The handback
#    file name + datetime get filled in, the target file (a.md) becomes a
#    hyperlinked "Latest Target File" entry, and the columns get re-sized.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 15:18:43"

$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 15:18:43"

# Rebuild the hyperlinks collection in reading order (A2, I2, A3, I3) so the
# relationship ids line up the way a full Excel save would emit them.
$dede.Hyperlinks.Delete()

$dede.Hyperlinks.Add(
    $dede.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a72fcfa1e9d351cfe847c3a15ce311c9913c1823/e2e/a.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "a.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a72fcfa1e9d351cfe847c3a15ce311c9913c1823/e2e/a.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "a.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a72fcfa1e9d351cfe847c3a15ce311c9913c1823/e2e/b.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "b.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a72fcfa1e9d351cfe847c3a15ce311c9913c1823/e2e/a.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "a.md"
)

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
